# This script updates the state-transition probability matrix on Sheet1.
# Several rows that previously had all-zero transition probabilities now
# contain the computed probabilities (reflecting additional simulated
# games / updated simulate-game logic).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4444444444444444
$ws.Range("C2").Value = 0.1111111111111111
$ws.Range("P2").Value = 0.2222222222222222
$ws.Range("S2").Value = 0.2222222222222222

$ws.Range("P3").Value = 1

$ws.Range("S4").Value = 1

$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.1666666666666667
$ws.Range("S6").Value = 0.3333333333333333

$ws.Range("J7").Value = 0.2222222222222222
$ws.Range("Q7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.6666666666666666

$ws.Range("B8").Value = 0.07407407407407407
$ws.Range("J8").Value = 0.1481481481481481
$ws.Range("Q8").Value = 0.1851851851851852
$ws.Range("R8").Value = 0.1481481481481481
$ws.Range("S8").Value = 0.4444444444444444

$ws.Range("B9").Value = 0.2857142857142857
$ws.Range("F9").Value = 0.1428571428571428
$ws.Range("S9").Value = 0.5714285714285714

$ws.Range("B10").Value = 0.03125
$ws.Range("D10").Value = 0.03125
$ws.Range("F10").Value = 0.0625
$ws.Range("J10").Value = 0.0625
$ws.Range("Q10").Value = 0.15625
$ws.Range("R10").Value = 0.125
$ws.Range("S10").Value = 0.53125

$ws.Range("G11").Value = 0.15
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.55

$ws.Range("G12").Value = 0.5454545454545454
$ws.Range("J12").Value = 0.2727272727272727
$ws.Range("S12").Value = 0.1818181818181818

$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("O15").Value = 0.1666666666666667
$ws.Range("S15").Value = 0.3333333333333333

$ws.Range("H16").Value = 0.3333333333333333
$ws.Range("S16").Value = 0.6666666666666666

$ws.Range("H17").Value = 0.3333333333333333
$ws.Range("J17").Value = 0.1666666666666667
$ws.Range("K17").Value = 0.25
$ws.Range("O17").Value = 0.08333333333333333
$ws.Range("S17").Value = 0.1666666666666667

$ws.Range("H18").Value = 0.2222222222222222
$ws.Range("J18").Value = 0.2222222222222222
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("O18").Value = 0.1111111111111111
$ws.Range("S18").Value = 0.3333333333333333

$ws.Range("F19").Value = 0.01492537313432836
$ws.Range("H19").Value = 0.2985074626865671
$ws.Range("I19").Value = 0.1044776119402985
$ws.Range("J19").Value = 0.208955223880597
$ws.Range("K19").Value = 0.1641791044776119
$ws.Range("O19").Value = 0.01492537313432836
$ws.Range("S19").Value = 0.1940298507462687
